# Apply the change described by the diff:
# - Insert a new row before row 201 (shifting old rows 201-218 down to 202-219)
# - Populate the new row 201 with updated values (same as the former row 201,
#   except Fecha, Precio minimo/maximo/promedio and Precio $/Kg change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201; this pushes existing rows 201-218 down to 202-219
$ws.Rows("201:201").Insert()

# Fill the newly inserted row 201 with its data
$ws.Range("A201").Value = 5
$ws.Range("B201").Value = "Macroferia Regional de Talca"
$ws.Range("C201").Value = "Maule"
$ws.Range("D201").Value = 44461
$ws.Range("D201").NumberFormat = $ws.Range("D202").NumberFormat
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 100112043
$ws.Range("G201").Value = "Pepino ensalada"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 400
$ws.Range("K201").Value = 16000
$ws.Range("L201").Value = 16000
$ws.Range("M201").Value = 16000
$ws.Range("N201").Value = "`$/caja 60 unidades"
$ws.Range("O201").Value = "Región de Arica y Parinacota"
$ws.Range("P201").Value = 267
$ws.Range("Q201").Value = 60
$ws.Range("R201").Value = "Hortaliza"
